# Apply conversion of the weekly schedule sheet:
#  - Row 1 header (B1:F1): numeric weekday numbers -> Portuguese weekday names
#  - Column A (A2:A13): numeric period index -> class start times
#  - A few subject cells shifted one time-slot earlier (rows 5-12)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: weekday names ---
$ws.Range("B1").Value = "segunda"
$ws.Range("C1").Value = "terça"
$ws.Range("D1").Value = "quarta"
$ws.Range("E1").Value = "quinta"
$ws.Range("F1").Value = "sexta"

# --- Column A: time labels instead of period numbers ---
$ws.Range("A2").Value  = "7:00"
$ws.Range("A3").Value  = "7:50"
$ws.Range("A4").Value  = "8:40"
$ws.Range("A5").Value  = "9:30"
$ws.Range("A6").Value  = "10:40"
$ws.Range("A7").Value  = "11:30"
$ws.Range("A8").Value  = "13:00"
$ws.Range("A9").Value  = "13:50"
$ws.Range("A10").Value = "14:40"
$ws.Range("A11").Value = "15:30"
$ws.Range("A12").Value = "16:40"
$ws.Range("A13").Value = "17:30"

# --- Subject cells that moved to a different time slot ---
$ws.Range("F5").Value  = "EAP"
$ws.Range("C6").Value  = "EAP"
$ws.Range("B7").Value  = "-"
$ws.Range("E7").Value  = "Desenho Técnico"
$ws.Range("B8").Value  = "Circuitos Elétricos 2"
$ws.Range("B10").Value = "-"
$ws.Range("C12").Value = "-"
$ws.Range("D12").Value = "-"
